$d = $word.ActiveDocument

# 1. "higher educational institutions" -> "higher education institutions"
#    (delete the "al" that turns "education" into "educational")
$d.Content.Find.Execute("higher educational institutions", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "higher education institutions", 2)

# 2. Move the auto "_GoBack" bookmark from just before "higher education institutions"
#    to the point between "...for some communities " and "and are rare in others."
#    (this mirrors Word re-stamping _GoBack at the location of the most recent edit)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$target = $d.Content.Duplicate
$target.Find.Execute("for some communities and are rare in others.", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)
$markPoint = $d.Range($target.Start + 21, $target.Start + 21)
$d.Bookmarks.Add("_GoBack", $markPoint)
